# Update view/attendance counts on the "展览" and "全部类型" sheets to
# reflect regenerated stats (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 7497
$wsExhibition.Range("F7").Value = 4091
$wsExhibition.Range("F9").Value = 574

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7497
$wsAll.Range("F9").Value = 4091
$wsAll.Range("F11").Value = 574
